# Update the division problems in the document to match the new output.
# Note: "76÷6=" is replaced before "79÷3=" is turned into "76÷6=" to avoid
# a collision where the newly-written text would get caught by a later
# Find/Replace operation.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "60÷3=" "87÷8="
Replace-Text "22÷9=" "13÷5="
Replace-Text "66÷2=" "53÷9="
Replace-Text "75÷8=" "97÷6="
Replace-Text "34÷6=" "28÷6="
Replace-Text "69÷6=" "57÷8="
Replace-Text "39÷6=" "50÷9="
Replace-Text "64÷3=" "37÷6="
Replace-Text "33÷5=" "62÷8="
Replace-Text "61÷3=" "44÷7="
Replace-Text "58÷3=" "45÷3="
Replace-Text "76÷6=" "98÷2="
Replace-Text "79÷3=" "76÷6="
Replace-Text "31÷8=" "47÷6="
Replace-Text "28÷5=" "91÷3="
Replace-Text "75÷7=" "56÷2="
Replace-Text "23÷7=" "35÷6="
Replace-Text "95÷6=" "91÷2="
Replace-Text "22÷3=" "50÷4="
Replace-Text "30÷6=" "51÷8="
Replace-Text "71÷2=" "49÷7="
Replace-Text "37÷3=" "67÷9="
Replace-Text "69÷2=" "36÷6="
Replace-Text "43÷5=" "38÷7="
Replace-Text "54÷8=" "30÷3="
